$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: insert a run of `text` at the end of the target paragraph,
# using the exact same run-properties (rFonts/color/sz/szCs/shd) as
# an existing "clean" run in the document (one with no rsid
# attributes on its <w:r>), while keeping it as its own separate
# <w:r> element (not merged into the neighbouring run).
#
# Strategy: FormattedText assignment into the document preserves
# formatting perfectly and never merges with neighbours. But editing
# the *text* of a freshly inserted run (Range.Text = ...) causes the
# engine to merge it back into the preceding run. To avoid that we
# first stage the exact wanted text+formatting in a disposable
# scratch paragraph appended at the very end of the document, read
# back its FormattedText (now carrying the right text AND right
# formatting), insert that into the real target location (clean,
# un-merged), and finally delete the scratch paragraph.
# ------------------------------------------------------------------

# Locate a template run already carrying the desired formatting
# (rFonts cstheme=minorHAnsi, color 24292F, sz/szCs 24, shd clear/FFFFFF)
# and with a plain <w:r> (no rsid attributes), so the copies we make
# are equally "clean".
$tmplRange = $d.Content.Duplicate
$null = $tmplRange.Find.Execute("The the master PC contains ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fmt = $tmplRange.FormattedText

function Append-FormattedRun($targetParaEnd, $text) {
    # Add disposable scratch paragraph at the end of the document.
    $scratchPara = $d.Paragraphs.Add()
    $scratchStart = $scratchPara.Range.Start

    # Stamp it with the template formatting.
    $scratchInsPt = $d.Range($scratchStart, $scratchStart)
    $scratchInsPt.FormattedText = $fmt

    # Replace its text with the text we actually want (merging here
    # is harmless - this paragraph is thrown away afterwards).
    $scratchParaNow = $d.Paragraphs.Last
    $scratchEnd = $scratchParaNow.Range.End
    $scratchRunRange = $d.Range($scratchStart, $scratchEnd - 1)
    $scratchRunRange.Text = $text

    # Capture the now correctly-formatted-and-texted run.
    $scratchParaFinal = $d.Paragraphs.Last
    $finalScratchEnd = $scratchParaFinal.Range.End
    $finalScratchRunRange = $d.Range($scratchStart, $finalScratchEnd - 1)
    $finalFmt = $finalScratchRunRange.FormattedText

    # Insert it, as a clean standalone run, right before the target
    # paragraph's end-of-paragraph mark.
    $insPt = $d.Range($targetParaEnd - 1, $targetParaEnd - 1)
    $insPt.FormattedText = $finalFmt

    # Remove the scratch paragraph again.
    $scratchParaToDelete = $d.Paragraphs.Last
    $scratchParaToDelete.Range.Delete()
}

# The paragraph we are appending the new sentences to is the very
# last paragraph of the document body.
$targetPara = $d.Paragraphs.Last
$targetParaEnd = $targetPara.Range.End
Append-FormattedRun $targetParaEnd "It is important to be"

$targetPara = $d.Paragraphs.Last
$targetParaEnd = $targetPara.Range.End
Append-FormattedRun $targetParaEnd " sure that the audio output is correctly configured on the client devices. Check that the volume is turned up and that the correct audio device is selected."

$targetPara = $d.Paragraphs.Last
$targetParaEnd = $targetPara.Range.End
Append-FormattedRun $targetParaEnd " "

$targetPara = $d.Paragraphs.Last
$targetParaEnd = $targetPara.Range.End
Append-FormattedRun $targetParaEnd "Check that the IP addresses and port numbers specified in the `"multiudpsink`" element of the master device pipeline match the `"address`" and `"port`" options specified in the client device pipelines."

Write-Output "done"
